$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "check" column (B) for completed checklist items
$ws.Range("B2").Value = "check"
$ws.Range("B3").Value = "check"
$ws.Range("B4").Value = "check"
$ws.Range("B5").Value = "check"
$ws.Range("B6").Value = "check"
$ws.Range("B7").Value = "check"
$ws.Range("B9").Value = "check"
$ws.Range("B10").Value = "check"
$ws.Range("B11").Value = "check"
$ws.Range("B13").Value = "check"

# Add note about which pin was actually used
$ws.Range("C13").Value = "used PB6"

# Update the active selection to reflect where the user left off
$ws.Range("B11").Select() | Out-Null
